$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.634.60"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "3.473.09"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.96%  "
$ws.Range("D7").Value = "3.468.92"
$ws.Range("E7").Value = "  -2.21%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("E13").Value = "  -3.85%  "
$ws.Range("D14").Value = "4.061.93"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.60%  "
$ws.Range("D16").Value = "3.480.59"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "66.709.26"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "437.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E23").Value = "  -5.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "3.613.84"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("E27").Value = "  -7.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.36"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.73%  "
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("E31").Value = "  -6.23%  "
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.30"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").Value = "3.465.68"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.17%  "
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "175.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0881"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("E43").Value = "  -10.98%  "
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.890"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.88"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.43"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.45%  "
$ws.Range("E50").Value = "  -8.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.977"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.46%  "
